$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AF2").Value = 7
$ws.Range("AI2").Value = 1

# Row 3
$ws.Range("D3").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AC3").Value = 7
$ws.Range("AF3").Value = 7

# Row 4
$ws.Range("D4").Value = "2024-07-08T01:58:00.000Z"

# Row 5
$ws.Range("D5").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AC5").Value = 6
$ws.Range("AF5").Value = 6.5

# Row 6
$ws.Range("D6").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AF6").Value = 7
$ws.Range("AI6").Value = 1

# Row 7
$ws.Range("D7").Value = "2024-07-08T01:58:00.000Z"

# Row 8
$ws.Range("D8").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AC8").Value = 6
$ws.Range("AF8").Value = 6.5

# Row 9
$ws.Range("D9").Value = "2024-07-08T01:58:00.000Z"

# Row 10
$ws.Range("D10").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AC10").Value = 5
$ws.Range("AF10").Value = 6.5

# Row 11
$ws.Range("D11").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AC11").Value = 7
$ws.Range("AF11").Value = 7

# Row 12
$ws.Range("D12").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AC12").Value = 7
$ws.Range("AF12").Value = 7

# Row 13
$ws.Range("D13").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AC13").Value = 7
$ws.Range("AF13").Value = 7

# Row 14
$ws.Range("D14").Value = "2024-07-08T01:58:00.000Z"
$ws.Range("AC14").Value = 7
$ws.Range("AF14").Value = 7
